# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reorder the "Periodo Mora" rows (16-20) into chronological order
# (2109, 2110, 2111, 2112, 2201) and refresh the "Valor Mora" (F) /
# "Salario Basico" (G) figures for the updated EC database.
$ws.Range("E16").Value = "2109"
$ws.Range("F16").Value = 36341
$ws.Range("G16").Value = 877803

$ws.Range("E17").Value = "2110"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 877803

$ws.Range("E18").Value = "2111"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 877803

$ws.Range("E19").Value = "2112"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 877803

$ws.Range("E20").Value = "2201"
$ws.Range("F20").Value = 29260
$ws.Range("G20").Value = 877803
